$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell values (swap pairs as per diff)
$ws.Range("C8").Value = 3
$ws.Range("E8").Value = 1

$ws.Range("C14").Value = 8
$ws.Range("E14").Value = 3
$ws.Range("I14").Value = 8
$ws.Range("K14").Value = 6

# Update the active selection on the sheet
$ws.Range("C15").Select()
